$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row 2's cell formatting into row 3 (skip column F, which has no
#     cell in row 2) so the new row matches the existing "Log" formatting. ---
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A3:E3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("G2:H2").Copy() | Out-Null
$ws.Range("G3:H3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- The existing 10:00-14:00 (4h) block on row 2 is split into two
#     separate, shorter sessions: 10:00-11:45 and 13:30-15:00. ---

# Row 2: shorten the End time from 14:00 to 11:45.
$ws.Range("C2").Value = 0.48958333333333331   # 11:45 AM

# Row 3: brand-new entry for the second (afternoon) work session.
$ws.Range("A3").Formula = "=DATE(2025,4,27)"
$ws.Range("B3").Value = 0.5625                # 1:30 PM
$ws.Range("C3").Value = 0.625                 # 3:00 PM
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = 12
$ws.Range("G3").Formula = "=(C3-B3)*24"
$ws.Range("H3").Formula = "=WEEKNUM(A3)"

# Recalculate formulas (Time column, WEEKNUM, pivot table output, ...).
$wb.Application.Calculate()

# Refresh the pivot table/cache so it reflects the new source rows.
$pt = $ws.PivotTables(1)
$pt.RefreshTable() | Out-Null

# Move the active selection, matching the author's last cursor position.
$ws.Range("C7").Select() | Out-Null
